# Weekly fruit/vegetable price update: a new weekly record is inserted
# at row 59 (pushing the existing historical rows down by one), matching
# the "Fruta / hortaliza, semanal" ingestion pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 59, shifting rows 59:92
# down to 60:93 (dimension grows from A1:R92 to A1:R93).
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with this week's record.
$ws.Cells.Item(59, 1).Value = 4
$ws.Cells.Item(59, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(59, 3).Value = "Los Lagos"
$ws.Cells.Item(59, 4).Value = 44582
$ws.Cells.Item(59, 5).Value = 10
$ws.Cells.Item(59, 6).Value = 100112052
$ws.Cells.Item(59, 7).Value = "Albahaca"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 140
$ws.Cells.Item(59, 11).Value = 6000
$ws.Cells.Item(59, 12).Value = 6000
$ws.Cells.Item(59, 13).Value = 6000
$ws.Cells.Item(59, 14).Value = "$/docena de matas"
$ws.Cells.Item(59, 15).Value = "Región Metropolitana"
$ws.Cells.Item(59, 16).Value = 1000
$ws.Cells.Item(59, 17).Value = 6
$ws.Cells.Item(59, 18).Value = "Hortaliza"
